$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "Q8"
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B2").Value = -1.123305982541588
$ws.Range("C2").Value = -0.02330598254158778
$ws.Range("D2").Value = 0.1766940174584122
$ws.Range("E2").Value = -0.2233059825415878
$ws.Range("F2").Value = 0.07669401745841223
$ws.Range("G2").Value = -0.3233059825415878
$ws.Range("H2").Value = 0.4766940174584122
$ws.Range("B3").Value = 0.043791487406589
$ws.Range("C3").Value = 0.243791487406589
$ws.Range("D3").Value = -0.156208512593411
$ws.Range("E3").Value = 0.143791487406589
$ws.Range("F3").Value = -0.256208512593411
$ws.Range("G3").Value = 0.543791487406589
$ws.Range("B4").Value = 0.125432791727192
$ws.Range("C4").Value = -0.2745672082728081
$ws.Range("D4").Value = 0.02543279172719195
$ws.Range("E4").Value = -0.3745672082728081
$ws.Range("F4").Value = 0.425432791727192
$ws.Range("G4").Value = 0.125432791727192
$ws.Range("H4").Value = 0.02543279172719195
$ws.Range("I4").Value = -0.3745672082728081
$ws.Range("J4").Value = -0.7745672082728081
$ws.Range("B5").Value = -0.3015462512541777
$ws.Range("C5").Value = -0.0015462512541777
$ws.Range("D5").Value = -0.4015462512541777
$ws.Range("E5").Value = 0.3984537487458223
$ws.Range("F5").Value = 0.0984537487458223
$ws.Range("G5").Value = -0.0015462512541777
$ws.Range("H5").Value = -0.4015462512541777
$ws.Range("I5").Value = -0.8015462512541778
$ws.Range("B6").Value = 0.271863973860668
$ws.Range("C6").Value = -0.128136026139332
$ws.Range("D6").Value = 0.671863973860668
$ws.Range("E6").Value = 0.3718639738606681
$ws.Range("F6").Value = 0.271863973860668
$ws.Range("G6").Value = -0.128136026139332
$ws.Range("H6").Value = -0.528136026139332
$ws.Range("B7").Value = -0.1414808197323844
$ws.Range("C7").Value = 0.6585191802676156
$ws.Range("D7").Value = 0.3585191802676156
$ws.Range("E7").Value = 0.2585191802676156
$ws.Range("F7").Value = -0.1414808197323844
$ws.Range("G7").Value = -0.5414808197323844
$ws.Range("B8").Value = 0.4029184009508193
$ws.Range("C8").Value = 0.1029184009508193
$ws.Range("D8").Value = 0.002918400950819272
$ws.Range("E8").Value = -0.3970815990491807
$ws.Range("F8").Value = -0.7970815990491807
$ws.Range("G8").Value = -0.7970815990491807
$ws.Range("H8").Value = 0.9029184009508193
$ws.Range("I8").Value = -0.4970815990491807
$ws.Range("B9").Value = 0.1966784556707163
$ws.Range("C9").Value = 0.09667845567071635
$ws.Range("D9").Value = -0.3033215443292837
$ws.Range("E9").Value = -0.7033215443292837
$ws.Range("F9").Value = -0.7033215443292837
$ws.Range("G9").Value = 0.9966784556707163
$ws.Range("H9").Value = -0.4033215443292836
$ws.Range("B10").Value = 0.05000130330345419
$ws.Range("C10").Value = -0.3499986966965458
$ws.Range("D10").Value = -0.7499986966965458
$ws.Range("E10").Value = -0.7499986966965458
$ws.Range("F10").Value = 0.9500013033034542
$ws.Range("G10").Value = -0.4499986966965458
$ws.Range("B11").Value = -0.4206125572440377
$ws.Range("C11").Value = -0.8206125572440377
$ws.Range("D11").Value = -0.8206125572440377
$ws.Range("E11").Value = 0.8793874427559624
$ws.Range("F11").Value = -0.5206125572440377
$ws.Range("B12").Value = -0.7551200626749693
$ws.Range("C12").Value = -0.7551200626749693
$ws.Range("D12").Value = 0.9448799373250307
$ws.Range("E12").Value = -0.4551200626749692
$ws.Range("B13").Value = -0.7352230408639261
$ws.Range("C13").Value = 0.964776959136074
$ws.Range("D13").Value = -0.435223040863926
$ws.Range("B14").Value = 1.027145135284297
$ws.Range("C14").Value = -0.3728548647157029
$ws.Range("B15").Value = -0.5726054543893956
